$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run with updated TPM data: the "FAPs" sending-cluster rows (currently
# rows 2-3) get pushed down to rows 4-5, with their specificity columns
# (I, J, P, Q, R, S, T) recomputed now that a second sending cluster
# ("ECs") also sends Ccl28->Ccr10 signal. New "ECs" rows take rows 2-3.

# Rows 4-5: original "FAPs" rows, recalculated specificities.
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Ccl28"
$ws.Cells.Item(4,3).Value = "Ccr10"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1070356666666667
$ws.Cells.Item(4,8).Value = 0.321107
$ws.Cells.Item(4,9).Value = 0.4505354815264836
$ws.Cells.Item(4,10).Value = 0.4505354815264836
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 1.036595333333333
$ws.Cells.Item(4,14).Value = 3.109786
$ws.Cells.Item(4,15).Value = 0.393072250513715
$ws.Cells.Item(4,16).Value = 0.3930722505137151
$ws.Cells.Item(4,17).Value = 0.1109526725668889
$ws.Cells.Item(4,18).Value = 0.998574053102
$ws.Cells.Item(4,19).Value = 0.1770929956598952
$ws.Cells.Item(4,20).Value = 0.1770929956598952

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ccl28"
$ws.Cells.Item(5,3).Value = "Ccr10"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.1070356666666667
$ws.Cells.Item(5,8).Value = 0.321107
$ws.Cells.Item(5,9).Value = 0.4505354815264836
$ws.Cells.Item(5,10).Value = 0.4505354815264836
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.600567
$ws.Cells.Item(5,14).Value = 4.801701
$ws.Cells.Item(5,15).Value = 0.6069277494862849
$ws.Cells.Item(5,16).Value = 0.6069277494862849
$ws.Cells.Item(5,17).Value = 0.1713177558896666
$ws.Cells.Item(5,18).Value = 1.541859803007
$ws.Cells.Item(5,19).Value = 0.2734424858665884
$ws.Cells.Item(5,20).Value = 0.2734424858665884

# Rows 2-3: new "ECs" sending-cluster rows.
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ccl28"
$ws.Cells.Item(2,3).Value = "Ccr10"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.1305386666666667
$ws.Cells.Item(2,8).Value = 0.391616
$ws.Cells.Item(2,9).Value = 0.5494645184735164
$ws.Cells.Item(2,10).Value = 0.5494645184735164
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.036595333333333
$ws.Cells.Item(2,14).Value = 3.109786
$ws.Cells.Item(2,15).Value = 0.393072250513715
$ws.Cells.Item(2,16).Value = 0.3930722505137151
$ws.Cells.Item(2,17).Value = 0.1353157726862222
$ws.Cells.Item(2,18).Value = 1.217841954176
$ws.Cells.Item(2,19).Value = 0.2159792548538198
$ws.Cells.Item(2,20).Value = 0.2159792548538199

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ccl28"
$ws.Cells.Item(3,3).Value = "Ccr10"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.1305386666666667
$ws.Cells.Item(3,8).Value = 0.391616
$ws.Cells.Item(3,9).Value = 0.5494645184735164
$ws.Cells.Item(3,10).Value = 0.5494645184735164
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.600567
$ws.Cells.Item(3,14).Value = 4.801701
$ws.Cells.Item(3,15).Value = 0.6069277494862849
$ws.Cells.Item(3,16).Value = 0.6069277494862849
$ws.Cells.Item(3,17).Value = 0.2089358820906666
$ws.Cells.Item(3,18).Value = 1.880422938816
$ws.Cells.Item(3,19).Value = 0.3334852636196965
$ws.Cells.Item(3,20).Value = 0.3334852636196965
